$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.045.21"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.906.65"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4802"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2983"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06654"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "101.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +19.45%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.896.88"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07657"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.185"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6681"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "308.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +25.29%  "
$ws.Range("D17").Value = "31.004.21"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.58%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007611"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").Value = "2.144.23"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.243"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.261"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.008"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.91%  "
$ws.Range("E29").Value = "  +9.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.363"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.205"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.034"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05116"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7566"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.176"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.761"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02012"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.725"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.081"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8887"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "109.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.64%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4256"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.44%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.749"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "68.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.428"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.261"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1241"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.440"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05680"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.83%  "
